$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: error code 412501 / EntityAlreadyExists description
$ws.Range("A3").Value = 412501
$ws.Range("B3").Value = "EntityAlreadyExists - The request was rejected because it attempted to create a resource that already exists."

# Match style of column A (left aligned) used by the existing row (A2 uses style index 3 / horizontal left)
$ws.Range("A3").HorizontalAlignment = -4131
